$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the H1 title.
#    Target shape:
#      <w:p>
#        <w:r/>
#        <w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r>
#        <w:r><w:t>: Explore gameplay, ... challenge your luck.</w:t></w:r>
#      </w:p>
# ---------------------------------------------------------------------------

$titlePara = $d.Paragraphs(1)
$gameplayHeadingPara = $d.Paragraphs(2)

# Insert a brand-new empty paragraph right before "Gameplay Features" (i.e.
# immediately after the title) and strip the inherited heading style so it
# becomes a normal body paragraph.
$gameplayHeadingPara.Range.InsertParagraphBefore() | Out-Null
$metaPara = $d.Paragraphs(2)
$metaPara.Range.Style = "Normal"

$metaRange = $metaPara.Range
$metaRange.InsertAfter("Meta description: Explore gameplay, symbols, winning potential, RTP, and bonus features in our review of Black Bull. Play for free and challenge your luck.")

# Bold just the "Meta description" lead-in (16 characters).
$metaStart = $d.Paragraphs(2).Range.Start
$boldRange = $d.Range($metaStart, $metaStart + 16)
$boldRange.Font.Bold = 1

# Leave a leading empty run, matching the pattern used by every other body
# paragraph in this document (each one opens with an empty <w:r/>).
$leadRange = $d.Range($metaStart, $metaStart)
$leadRange.InsertBefore("")

# ---------------------------------------------------------------------------
# 2) Remove the duplicated bold title paragraph near the end of the document
#    and rewrite the italic paragraph that follows it with the new image
#    prompt text (keeping its italic formatting intact).
# ---------------------------------------------------------------------------

$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($count)
$secondLastPara = $d.Paragraphs($count - 1)

$oldTitleText = "Play Black Bull for Free: Review of Classic Slot"
$oldMetaSentence = "Explore gameplay, symbols, winning potential, RTP, and bonus features in our review of Black Bull. Play for free and challenge your luck."

if ($secondLastPara.Range.Text.TrimEnd([char]13, [char]7) -eq $oldTitleText -and `
    $lastPara.Range.Text.TrimEnd([char]13, [char]7) -eq $oldMetaSentence) {
    $secondLastPara.Range.Delete()
}

$count = $d.Paragraphs.Count
$finalPara = $d.Paragraphs($count)
$finalRange = $finalPara.Range
$finalRange.MoveEnd(1, -1) | Out-Null

$newPromptText = 'Please create a cartoon-style feature image of a happy Maya warrior with glasses for the online slot game "Black Bull". The image should showcase the warrior in a victorious pose, surrounded by the symbols from the game, including the black bull, wolves, eagles, geckos, and playing card symbols. The colors should be bright and vibrant, and the image should be eye-catching to grab the attention of players. Use creative license to add any fun elements or features that will enhance the image and make it stand out.'

$finalRange.Text = $newPromptText

Write-Output "Paragraph count: $($d.Paragraphs.Count)"
Write-Output "Meta paragraph: $($d.Paragraphs(2).Range.Text)"
Write-Output "Final paragraph: $($d.Paragraphs($d.Paragraphs.Count).Range.Text)"
